# Update the "concise_ms" marking-scheme pattern:
#   Row 11 = per-question marks for Right / Wrong / Not-Attempted
#   Row 12 = resulting Total score and "score/max" summary string
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marking scheme: +5 for a right answer, -1.2 for a wrong one
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Totals recomputed from the new marking scheme
# (13 right * 5 = 65, 0 wrong * -1.2 = -0, max = 28 * 5 = 140)
$ws.Range("B12").Value = 65
$ws.Range("C12").Value = -0

$ws.Range("E12").Value = "65.0/140"
